$wb = $excel.ActiveWorkbook

# --- Reference sheet used purely to copy cell formatting (style) ---
$fmtSrc = $wb.Worksheets.Item("2021-Q4")

# --- Step 1: remove the existing "总计" sheet; we'll recreate it after the
#     new quarter sheet so the engine hands out sheetId 6 / 7 in the order
#     the target workbook expects. ---
$zjOld = $wb.Worksheets.Item("总计")
$zjOld.Delete()

# --- Step 2: insert the new "2022-Q1" sheet right after the last
#     remaining sheet (currently "2021-Q4"), matching its position in the
#     tab order. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($null, $lastSheet)
$q1.Name = "2022-Q1"
# Header row, formatted like the other quarter sheets
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Data rows
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'161725"
$q1.Range("C2").Value = "招商中证白酒指数"
$q1.Range("D2").Value = "'688.84"
$q1.Range("E2").Value = "'94.34"
$q1.Range("F2").Value = "'2.59"
$q1.Range("G2").Value = "'17.8410"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'001705"
$q1.Range("C3").Value = "泓德战略转型股票"
$q1.Range("D3").Value = "'22.82"
$q1.Range("E3").Value = "'93.94"
$q1.Range("F3").Value = "'3.85"
$q1.Range("G3").Value = "'0.8786"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'006926"
$q1.Range("C4").Value = "长城量化精选股票A"
$q1.Range("D4").Value = "'4.97"
$q1.Range("E4").Value = "'90.86"
$q1.Range("F4").Value = "'7.90"
$q1.Range("G4").Value = "'0.3926"
$q1.Range("H4").Value = 7

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'011013"
$q1.Range("C5").Value = "长城消费30股票型证券投资基金A"
$q1.Range("D5").Value = "'4.93"
$q1.Range("E5").Value = "'90.83"
$q1.Range("F5").Value = "'3.27"
$q1.Range("G5").Value = "'0.1612"
$q1.Range("H5").Value = 9

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "'011463"
$q1.Range("C6").Value = "长城量化精选股票C"
$q1.Range("D6").Value = "'1.09"
$q1.Range("E6").Value = "'90.86"
$q1.Range("F6").Value = "'7.90"
$q1.Range("G6").Value = "'0.0861"
$q1.Range("H6").Value = 7

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "'001695"
$q1.Range("C7").Value = "泓德泓业灵活配置混合"
$q1.Range("D7").Value = "'1.16"
$q1.Range("E7").Value = "'91.88"
$q1.Range("F7").Value = "'3.52"
$q1.Range("G7").Value = "'0.0408"
$q1.Range("H7").Value = 10

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "'200016"
$q1.Range("C8").Value = "长城稳健成长灵活配置混合"
$q1.Range("D8").Value = "'0.83"
$q1.Range("E8").Value = "'78.39"
$q1.Range("F8").Value = "'4.59"
$q1.Range("G8").Value = "'0.0381"
$q1.Range("H8").Value = 4

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "'009015"
$q1.Range("C9").Value = "泓德睿享一年持有期混合A"
$q1.Range("D9").Value = "'3.41"
$q1.Range("E9").Value = "'24.56"
$q1.Range("F9").Value = "'0.92"
$q1.Range("G9").Value = "'0.0314"
$q1.Range("H9").Value = 7

$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "'002159"
$q1.Range("C10").Value = "东吴国企改革主题灵活配置混合"
$q1.Range("D10").Value = "'0.24"
$q1.Range("E10").Value = "'90.26"
$q1.Range("F10").Value = "'7.24"
$q1.Range("G10").Value = "'0.0174"
$q1.Range("H10").Value = 8

$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "'011014"
$q1.Range("C11").Value = "长城消费30股票型证券投资基金C"
$q1.Range("D11").Value = "'0.40"
$q1.Range("E11").Value = "'90.83"
$q1.Range("F11").Value = "'3.27"
$q1.Range("G11").Value = "'0.0131"
$q1.Range("H11").Value = 9

$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "'009527"
$q1.Range("C12").Value = "浙商汇金新兴消费灵活配置混合"
$q1.Range("D12").Value = "'0.34"
$q1.Range("E12").Value = "'73.34"
$q1.Range("F12").Value = "'3.62"
$q1.Range("G12").Value = "'0.0123"
$q1.Range("H12").Value = 9

$q1.Range("A13").Value = 11
$q1.Range("B13").Value = "'004805"
$q1.Range("C13").Value = "长信消费精选行业量化股票"
$q1.Range("D13").Value = "'0.11"
$q1.Range("E13").Value = "'89.89"
$q1.Range("F13").Value = "'7.76"
$q1.Range("G13").Value = "'0.0085"
$q1.Range("H13").Value = 4

$q1.Range("A14").Value = 12
$q1.Range("B14").Value = "'009016"
$q1.Range("C14").Value = "泓德睿享一年持有期混合C"
$q1.Range("D14").Value = "'0.07"
$q1.Range("E14").Value = "'24.56"
$q1.Range("F14").Value = "'0.92"
$q1.Range("G14").Value = "'0.0006"
$q1.Range("H14").Value = 7

# Copy column-A style (bold/border/center) from the reference sheet for A2:A14
$fmtSrc.Range("A2").Copy()
$q1.Range("A2:A14").PasteSpecial(-4122)

# Clear the incidental "Text" number-format the engine applies when a
# numeric-looking string is forced to text, so these cells end up with no
# explicit style (matching the other quarter sheets).
$q1.Range("B2:B14").ClearFormats()
$q1.Range("D2:G14").ClearFormats()

# --- Step 3: re-create the "总计" sheet after the new "2022-Q1" sheet so
#     it gets sheetId 7 (matching the target workbook). ---
$zj = $wb.Worksheets.Add($null, $q1)
$zj.Name = "总计"

$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"
$fmtSrc.Range("B1:D1").Copy()
$zj.Range("B1:D1").PasteSpecial(-4122)

# Data rows
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 13
$zj.Range("D2").Value = 19.52

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2021-Q4"
$zj.Range("C3").Value = 28
$zj.Range("D3").Value = 28.23

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q3"
$zj.Range("C4").Value = 66
$zj.Range("D4").Value = 47.73

$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q2"
$zj.Range("C5").Value = 40
$zj.Range("D5").Value = 20.86

$zj.Range("A6").Value = 4
$zj.Range("B6").Value = "2021-Q1"
$zj.Range("C6").Value = 26
$zj.Range("D6").Value = 4.79

$zj.Range("A7").Value = 5
$zj.Range("B7").Value = "2020-Q4"
$zj.Range("C7").Value = 24
$zj.Range("D7").Value = 16.4

# Copy column-A style (bold/border/center) from the reference sheet for A2:A7
$fmtSrc.Range("A2").Copy()
$zj.Range("A2:A7").PasteSpecial(-4122)
